$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab from "Scanner" to "Anatomy"
$ws.Name = "Anatomy"

# Delete the third data row (row 3) - the entire row so the table shrinks
$ws.Rows.Item(3).Delete()

# Re-mark the remaining used range (A1:F2) as "number stored as text" ignored,
# matching the shrunk ignoredErrors sqref (xlNumberAsText = 3)
$ws.Range("A1:F2").Errors.Item(3).Ignore = $true
